# Updates cell values in the active sheet to refresh the cryptos price list
# (values/percentages refreshed; two coin pairs re-ranked/swapped rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store numeric-looking text (e.g. '4.00',
# '0.999', '0.0000212') as text rather than coercing it to a Double and
# losing the exact printed form the source data uses.

$ws.Range('D2').Value = '91.870.62'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '3.182.12'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''241.48'
$ws.Range('E5').Value = '  +3.62%  '
$ws.Range('D6').Value = '''623.19'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = '''1.14'
$ws.Range('E7').Value = '  +5.81%  '
$ws.Range('D8').Value = '''0.374'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '3.181.89'
$ws.Range('E10').Value = '  +2.98%  '
$ws.Range('D11').Value = '''0.754'
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '''0.206'
$ws.Range('E12').Value = '  +4.30%  '
$ws.Range('D14').Value = '''35.80'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '91.701.71'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '3.748.12'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('D18').Value = '3.172.54'
$ws.Range('E18').Value = '  +2.51%  '
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('D20').Value = '''15.44'
$ws.Range('E20').Value = '  +10.16%  '
$ws.Range('B21').Value = 'PEPE'
$ws.Range('C21').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D21').Value = '''0.0000212'
$ws.Range('E21').Value = '  -5.32%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = '''5.95'
$ws.Range('E22').Value = '  +7.25%  '
$ws.Range('D23').Value = '''447.03'
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('D24').Value = '''9.28'
$ws.Range('E24').Value = '  +4.05%  '
$ws.Range('D25').Value = '''6.18'
$ws.Range('E25').Value = '  +8.63%  '
$ws.Range('D26').Value = '''89.56'
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('D27').Value = '''12.18'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').Value = '3.323.36'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '''0.135'
$ws.Range('E30').Value = '  +52.89%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').Value = '''0.172'
$ws.Range('E31').Value = '  +7.45%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.229'
$ws.Range('E32').Value = '  +16.38%  '
$ws.Range('D33').Value = '''9.48'
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('E34').Value = '  +9.28%  '
$ws.Range('D35').Value = '''7.78'
$ws.Range('E35').Value = '  +7.30%  '
$ws.Range('D36').Value = '''26.66'
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('D37').Value = '''517.22'
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').Value = '''0.889'
$ws.Range('E38').Value = '  -10.95%  '
$ws.Range('B39').Value = 'PancakeSwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D39').Value = '''1.94'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '''0.466'
$ws.Range('E40').Value = '  +14.83%  '
$ws.Range('B41').Value = 'MantraDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D41').Value = '''4.00'
$ws.Range('E41').Value = '  +16.47%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '''1.35'
$ws.Range('E42').Value = '  +5.14%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '''3.51'
$ws.Range('E43').Value = '  -8.37%  '
$ws.Range('D44').Value = '''22.18'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +4.05%  '
$ws.Range('D47').Value = '''1.94'
$ws.Range('E47').Value = '  +1.73%  '
$ws.Range('D48').Value = '''157.26'
$ws.Range('E48').Value = '  +3.05%  '
$ws.Range('D49').Value = '''1.40'
$ws.Range('E49').Value = '  +4.84%  '
$ws.Range('E50').Value = '  +14.57%  '
$ws.Range('E51').Value = '  +1.21%  '
